# Auto-generated edit script: updates crypto price/volume table
# to reflect the latest scrape (GitHub Actions bot run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.409.35"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "2.963.43"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "3.431.56"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").Value = "2.968.29"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.959"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").Value = "51.361.24"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +19.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.170"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("E29").Value = "  +9.37%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("E34").Value = "  -2.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.64%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("E40").Value = "  -3.66%  "

$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("E42").Value = "  +2.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.286"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("D48").Value = "2.027.98"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0338"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
